$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Add Panels"
$ws2 = $wb.Worksheets.Item(2)   # "Add Devices"

# ---------------------------------------------------------------------------
# Sheet 2 ("Add Devices"): add a new "Status" column (D) next to the existing
# Panel Name / Device Name / Device Type table, with True/False values that
# flag which device types are enabled, plus a "Status" header.
# ---------------------------------------------------------------------------

# Seed the shared-string table in the same order the authored workbook uses:
# True, False, Status (indices 37, 38, 39) - write the "True" rows first, the
# "False" rows second and the header last so the new strings land in that
# order.
$ws2.Range("D10").Value = "'True"
$ws2.Range("D11").Value = "'True"

$ws2.Range("D2").Value = "'False"
$ws2.Range("D3").Value = "'False"
$ws2.Range("D4").Value = "'False"
$ws2.Range("D5").Value = "'False"
$ws2.Range("D6").Value = "'False"
$ws2.Range("D7").Value = "'False"
$ws2.Range("D8").Value = "'False"
$ws2.Range("D9").Value = "'False"
$ws2.Range("D12").Value = "'False"

$ws2.Range("D1").Value = "Status"

# Give the new cells the same border/fill formatting as the rest of the
# table (column C) so the new column matches the existing look; the header
# cell picks up the same highlighted style as the other header cells.
$ws2.Range("C1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Range("C2:C12").Copy()
$ws2.Range("D2:D12").PasteSpecial(-4122)  # xlPasteFormats

# Re-assert the cell values (PasteSpecial only touched formatting, but make
# sure nothing was clobbered) and keep the text forced via the leading
# apostrophe so "True"/"False" are stored as text, not booleans.
$ws2.Range("D1").Value = "Status"
$ws2.Range("D2").Value = "'False"
$ws2.Range("D3").Value = "'False"
$ws2.Range("D4").Value = "'False"
$ws2.Range("D5").Value = "'False"
$ws2.Range("D6").Value = "'False"
$ws2.Range("D7").Value = "'False"
$ws2.Range("D8").Value = "'False"
$ws2.Range("D9").Value = "'False"
$ws2.Range("D10").Value = "'True"
$ws2.Range("D11").Value = "'True"
$ws2.Range("D12").Value = "'False"

# ---------------------------------------------------------------------------
# Selections: sheet1's active cell moves from H11 to H8, sheet2's moves from
# A21 to F5. Restore "Add Devices" as the active sheet afterwards, matching
# the workbook's original active-tab.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H8").Select()

$ws2.Activate()
$ws2.Range("F5").Select()
